{"js": "// Update the \"Heading 2\" (T\u00edtulo 2) paragraph style: increase the\n// space-before from 12pt (240 twips) to 24pt (480 twips) and pin\n// space-after to 0pt (0 twips), matching the updated style sheet.\nconst style = context.document.getStyles().getByNameOrNullObject(\"Heading 2\");\nstyle.load(\"nameLocal\");\nawait context.sync();\n\nif (style.isNullObject) {\n  throw new Error('Style \"Heading 2\" was not found in this document.');\n}\n\nconst pf = style.paragraphFormat;\npf.spaceBefore = 24; // points (was 12pt / 240 twips -> now 24pt / 480 twips)\npf.spaceAfter = 0;   // points (explicit 0, matches w:after=\"0\")\n\nawait context.sync();\n", "ps1": "# Update the \"Heading 2\" (T\u00edtulo 2) paragraph style: increase the\n# space-before from 12pt (240 twips) to 24pt (480 twips) and pin\n# space-after to 0pt (0 twips), matching the updated style sheet.\n$d = $word.ActiveDocument\n$style = $d.Styles(\"Heading 2\")\n$style.ParagraphFormat.SpaceBefore = 24\n$style.ParagraphFormat.SpaceAfter = 0\n"}
